$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New report row from DGS's 2021/09/29 report, appended right after the
# existing last data row (row 86 -> new row 87).
$newRow = 87

# Column A holds the report date as literal text (e.g. "2021/09/27" in the
# row above), not a real Excel date. Assigning a date-shaped string straight
# to .Value makes Excel auto-convert it to a date serial, which would change
# the cell's type/format. To avoid that, stage the text in a scratch cell
# that's forced to Text format, copy *values only* into the target cell
# (this preserves the target's pre-existing date-format style), then remove
# the scratch cell by shifting cells up so the sheet's used range doesn't
# grow beyond the new data.
$scratch = $ws.Cells.Item(1, 7)
$scratch.NumberFormat = "@"
$scratch.Value = "2021/09/29"

$dateCell = $ws.Cells.Item($newRow, 1)
$scratch.Copy()
$dateCell.PasteSpecial(-4163)
$scratch.Delete(-4159)

$ws.Cells.Item($newRow, 2).Value = 105.6
$ws.Cells.Item($newRow, 3).Value = 107.3
$ws.Cells.Item($newRow, 4).Value = 0.87
$ws.Cells.Item($newRow, 5).Value = 0.87

$ws.Range("A88").Select()
